$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1308.9
$ws.Range("I40").Value = 1238
$ws.Range("J40").Value = 1379.8
$ws.Range("K40").Value = 1238
$ws.Range("L40").Value = 1379.8
$ws.Range("M40").Value = -1063
$ws.Range("N40").Value = -1729.8
$ws.Range("H127").Value = 1263.1818
$ws.Range("I127").Value = 1059.8
$ws.Range("J127").Value = 1432.6666
$ws.Range("K127").Value = 3179.4
$ws.Range("L127").Value = 4297.9998
$ws.Range("M127").Value = 1780.6
$ws.Range("N127").Value = -14217.9998
$ws.Range("H128").Value = 46000
$ws.Range("J128").Value = 46000
$ws.Range("L128").Value = 46000
$ws.Range("N128").Value = -55960
$ws.Range("H138").Value = 136376.78
$ws.Range("I138").Value = 2182.8
$ws.Range("J138").Value = 180374.8
$ws.Range("K138").Value = 6548.400000000001
$ws.Range("L138").Value = 541124.3999999999
$ws.Range("M138").Value = -1408.400000000001
$ws.Range("N138").Value = -551404.3999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 18499.572
$ws.Range("J44").Value = 18499.572
$ws.Range("L44").Value = 18499.572
$ws.Range("N44").Value = -19475.572
$ws.Range("H55").Value = 29999
$ws.Range("J55").Value = 29999
$ws.Range("L55").Value = 29999
$ws.Range("N55").Value = -30629
$ws.Range("H76").Value = 100000
$ws.Range("J76").Value = 100000
$ws.Range("L76").Value = 100000
$ws.Range("N76").Value = -100676
$ws.Range("H79").Value = 100000
$ws.Range("J79").Value = 100000
$ws.Range("L79").Value = 100000
$ws.Range("N79").Value = -102340
$ws.Range("H94").Value = 30330
$ws.Range("J94").Value = 30330
$ws.Range("L94").Value = 30330
$ws.Range("N94").Value = -32132
$ws.Range("H132").Value = 4313.0884
$ws.Range("I132").Value = 3506.5715
$ws.Range("J132").Value = 5615.923
$ws.Range("K132").Value = 10519.7145
$ws.Range("L132").Value = 16847.769
$ws.Range("M132").Value = -7989.7145
$ws.Range("N132").Value = -21907.769

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 47657
$ws.Range("J76").Value = 47657
$ws.Range("L76").Value = 47657
$ws.Range("N76").Value = -48287
$ws.Range("H79").Value = 47657
$ws.Range("J79").Value = 47657
$ws.Range("L79").Value = 47657
$ws.Range("N79").Value = -49841
$ws.Range("H93").Value = 29525
$ws.Range("J93").Value = 29525
$ws.Range("L93").Value = 29525
$ws.Range("N93").Value = -33269
$ws.Range("H141").Value = 50080
$ws.Range("J141").Value = 50080
$ws.Range("L141").Value = 50080
$ws.Range("N141").Value = -60440

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3601.1333
$ws.Range("I134").Value = 2334.6667
$ws.Range("J134").Value = 5500.8335
$ws.Range("K134").Value = 7004.000100000001
$ws.Range("L134").Value = 16502.5005
$ws.Range("M134").Value = -4469.000100000001
$ws.Range("N134").Value = -21572.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2859671.5
$ws.Range("I4").Value = 5716329
$ws.Range("J4").Value = 3014.2856
$ws.Range("K4").Value = 17148987
$ws.Range("L4").Value = 9042.856800000001
$ws.Range("M4").Value = -17148875
$ws.Range("N4").Value = -9266.856800000001
$ws.Range("H113").Value = 1264.9
$ws.Range("I113").Value = 800
$ws.Range("J113").Value = 1419.8667
$ws.Range("K113").Value = 2400
$ws.Range("L113").Value = 4259.6001
$ws.Range("M113").Value = -230
$ws.Range("N113").Value = -8599.6001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2436.8
$ws.Range("I132").Value = 2010.3334
$ws.Range("J132").Value = 3076.5
$ws.Range("K132").Value = 6031.0002
$ws.Range("L132").Value = 9229.5
$ws.Range("M132").Value = -3501.0002
$ws.Range("N132").Value = -14289.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 9123.923000000001
$ws.Range("I22").Value = 833.6667
$ws.Range("J22").Value = 16229.857
$ws.Range("K22").Value = 833.6667
$ws.Range("L22").Value = 16229.857
$ws.Range("M22").Value = -538.6667
$ws.Range("N22").Value = -16819.857
$ws.Range("H27").Value = 9123.923000000001
$ws.Range("I27").Value = 833.6667
$ws.Range("J27").Value = 16229.857
$ws.Range("K27").Value = 833.6667
$ws.Range("L27").Value = 16229.857
$ws.Range("M27").Value = -726.6667
$ws.Range("N27").Value = -16443.857
$ws.Range("H59").Value = 47159
$ws.Range("J59").Value = 47159
$ws.Range("L59").Value = 47159
$ws.Range("N59").Value = -48467
$ws.Range("H68").Value = 3889.7585
$ws.Range("I68").Value = 2000
$ws.Range("J68").Value = 4609.6665
$ws.Range("K68").Value = 2000
$ws.Range("L68").Value = 4609.6665
$ws.Range("M68").Value = -1251
$ws.Range("N68").Value = -6107.6665
$ws.Range("H71").Value = 3889.7585
$ws.Range("I71").Value = 2000
$ws.Range("J71").Value = 4609.6665
$ws.Range("K71").Value = 10000
$ws.Range("L71").Value = 23048.3325
$ws.Range("M71").Value = -6256
$ws.Range("N71").Value = -30536.3325
$ws.Range("H75").Value = 66000
$ws.Range("J75").Value = 66000
$ws.Range("L75").Value = 66000
$ws.Range("N75").Value = -67872
$ws.Range("H78").Value = 66000
$ws.Range("J78").Value = 66000
$ws.Range("L78").Value = 198000
$ws.Range("N78").Value = -207360
$ws.Range("H105").Value = 75000
$ws.Range("J105").Value = 75000
$ws.Range("L105").Value = 75000
$ws.Range("N105").Value = -81988
$ws.Range("H132").Value = 4221.8647
$ws.Range("I132").Value = 3650.95
$ws.Range("K132").Value = 10952.85
$ws.Range("M132").Value = -8422.849999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 41727.9
$ws.Range("J46").Value = 41727.9
$ws.Range("L46").Value = 41727.9
$ws.Range("N46").Value = -42189.9
$ws.Range("H56").Value = 23905.6
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 23905.6
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 23905.6
$ws.Range("M56").ClearContents()
$ws.Range("N56").Value = -25333.6
$ws.Range("H132").Value = 3032589.2
$ws.Range("I132").Value = 2563.3103
$ws.Range("J132").Value = 6412233.5
$ws.Range("K132").Value = 7689.9309
$ws.Range("L132").Value = 19236700.5
$ws.Range("M132").Value = -5159.9309
$ws.Range("N132").Value = -19241760.5
$ws.Range("H134").Value = 41727.9
$ws.Range("J134").Value = 41727.9
$ws.Range("L134").Value = 125183.7
$ws.Range("N134").Value = -130253.7
